# Update the "想去人数" (F column) figures on the "展览" and "全部类型"
# sheets to reflect the latest scraped counts.

$wb = $excel.ActiveWorkbook

# Map of row number -> new value for column F
$updates = @{
    2  = 8361
    3  = 7846
    4  = 129
    6  = 38
    9  = 123
    14 = 1489
    19 = 127
}

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
